# Apply the cryptos price-list refresh described in the commit:
#   "Updated cryptos list on Sat May 13 08:58:41 UTC 2023 with GitHub Actions"
# Price (D) and Volume(1h) (E) columns are refreshed for nearly every row;
# rows 13/14 and 41/42 also swap their ranking order (Coin name + Link + Price + Volume).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.069.75'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +2.72%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.823.34'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +3.32%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.014'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +1.22%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.23'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +3.50%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.010'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.87%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4318'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.63%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3708'
$ws.Range("D8").ClearFormats()

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07290'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +3.46%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '2.159.87'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +24.14%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8731'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +4.73%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.42'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +6.09%  '

$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.662'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +3.78%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.430'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +3.76%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07006'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +3.24%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '81.32'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +2.93%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.016'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +1.12%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008951'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +4.15%  '

$ws.Range("E19").Value = '  +0.81%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.29'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +2.18%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '27.118.80'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +2.87%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.221'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +4.32%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.412.40'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +22.61%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.06'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.14%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.83'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.44%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.897'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +2.36%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.48'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +2.10%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.254'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +3.86%  '

$ws.Range("E29").Value = '  +13.78%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '115.26'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.81%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08969'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.67%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.180'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +6.93%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7507'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +3.60%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.457'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +3.12%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.825'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +2.76%  '

$ws.Range("E36").Value = '  +0.96%  '

$ws.Range("E37").Value = '  +5.31%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05256'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +2.92%  '

$ws.Range("E39").Value = '  +2.42%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5149'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +4.94%  '

$ws.Range("B41").Value = 'MXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.760'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +10.09%  '

$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1658'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +3.72%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.512'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +5.01%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.364'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +4.27%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '107.55'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +2.47%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.47'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +4.16%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.011'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.92%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4605'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +3.00%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.660'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +5.52%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06331'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +2.29%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.827'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +6.15%  '
